$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the red "// Clarify more where its used" review note, leaving
#    only the trailing line break run intact.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute(" //Clarify more where its used", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Delete()
}

# ---------------------------------------------------------------------------
# 2) Merge "To avoid them, we are hashing the passwords with " + "BCrypt" +
#    ", which is currently the industry " (previously split across runs by
#    a spell-check proofErr wrapper around "BCrypt") into a single run,
#    without touching the following "standard" run.
# ---------------------------------------------------------------------------
$bcryptOld = "To avoid them, we are hashing the passwords with BCrypt, which is currently the industry "
$r2 = $d.Content
$found2 = $r2.Find.Execute($bcryptOld, $false, $false, $false, $false, $false, $true, 1, $false, "BCRYPT-PLACEHOLDER-TOKEN", 2)
if ($found2) {
    $r2b = $d.Content
    $found2b = $r2b.Find.Execute("BCRYPT-PLACEHOLDER-TOKEN", $false, $false, $false, $false, $false, $true, 1, $false, $bcryptOld, 2)

    # The text-rebuild above also absorbs the following "standard" run into
    # the same run (since it shares identical formatting). Re-establish it
    # as its own run by round-tripping its font size.
    $r2c = $d.Content
    $found2c = $r2c.Find.Execute("standard", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2c) {
        $savedSize = $r2c.Font.Size
        $r2c.Font.Size = $savedSize + 1
        $r2c.Font.Size = $savedSize
    }
}

# ---------------------------------------------------------------------------
# 3) Merge "To prevent this, usage of 2" + "FA(" +
#    "2 factor authentication) is advised, as well as security checks
#    against weak passwords (top 10000 worst passwords)" into one run.
# ---------------------------------------------------------------------------
$twofaOld = "To prevent this, usage of 2FA(2 factor authentication) is advised, as well as security checks against weak passwords (top 10000 worst passwords)"
$r3 = $d.Content
$found3 = $r3.Find.Execute("To prevent this, usage of 2FA(2 factor authentication) is advised, as well as security checks against weak passwords (top 10000 worst passwords)", $false, $false, $false, $false, $false, $true, 1, $false, "TWOFA-PLACEHOLDER-TOKEN", 2)
if ($found3) {
    $r3b = $d.Content
    $found3b = $r3b.Find.Execute("TWOFA-PLACEHOLDER-TOKEN", $false, $false, $false, $false, $false, $true, 1, $false, $twofaOld, 2)
}

# ---------------------------------------------------------------------------
# 4) Merge "A10: " + "Server side" + " request forgery" into one run.
# ---------------------------------------------------------------------------
$a10Old = "A10: Server side request forgery"
$r4 = $d.Content
$found4 = $r4.Find.Execute("A10: Server side request forgery", $false, $false, $false, $false, $false, $true, 1, $false, "A10-PLACEHOLDER-TOKEN", 2)
if ($found4) {
    $r4b = $d.Content
    $found4b = $r4b.Find.Execute("A10-PLACEHOLDER-TOKEN", $false, $false, $false, $false, $false, $true, 1, $false, $a10Old, 2)
}

Write-Output "done: $found1 $found2 $found3 $found4"
